# Update Rizka - 26 Mei 2020
# Edits the "Admin - New Question" sheet: a couple of previously-duplicated
# "technicalCompetence/levelCompetence" cells are cleared and flagged red,
# a question's answer text is corrected, row 4 is refreshed from row 3's
# module/sub-module/question values, and the answer-key columns (I4/O4) are
# swapped from "uraian" to "kunciJawaban".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin - New Question")
$ws.Activate()

# Row 2: F2/G2 ("Lalala"/"Lalala") get cleared and highlighted red.
$ws.Range("F2").ClearContents()
$ws.Range("F2").Interior.Color = 255

$ws.Range("G2").ClearContents()
$ws.Range("G2").Interior.Color = 255

# Row 4: refresh tipeSoal/jobPosition/module/subModule from row 3's values,
# fix technicalCompetence, set levelCompetence to a number, move the answer
# from "uraian" (O4) to "kunciJawaban" (I4). Clear O4 before writing any new
# "Coba ya" text so the now-unused "Jawaban Uraian" shared-string slot is
# recycled for "Salah" first, matching the author's save order.
$ws.Range("A4").Value = "Benar / Salah"
$ws.Range("C4").Value = "SALES MANAGER"
$ws.Range("D4").Value = "NEOP RCCA"
$ws.Range("E4").Value = "Post-Test"

$ws.Range("O4").ClearContents()
$ws.Range("I4").Value = "Salah"

# Row 3: F3 gets corrected text; G3 gets cleared and highlighted red.
$ws.Range("F3").Value = "Coba ya"

$ws.Range("G3").ClearContents()
$ws.Range("G3").Interior.Color = 255

$ws.Range("F4").Value = "Coba ya"
$ws.Range("G4").Value = 12

# Restore the selection to where the author last left the cursor.
$ws.Range("E12").Select() | Out-Null
